$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3318.3684
$ws.Range("J64").Value = 3207.4666
$ws.Range("L64").Value = 3207.4666
$ws.Range("N64").Value = -3703.4666

$ws.Range("H67").Value = 3318.3684
$ws.Range("J67").Value = 3207.4666
$ws.Range("L67").Value = 3207.4666
$ws.Range("N67").Value = -4923.4666

$ws.Range("H107").Value = 247.07692
$ws.Range("I107").Value = 249.8
$ws.Range("J107").Value = 238
$ws.Range("K107").Value = 249.8
$ws.Range("L107").Value = 238
$ws.Range("M107").Value = 1670.2
$ws.Range("N107").Value = -4078

$ws.Range("H137").Value = 2719.111
$ws.Range("I137").Value = 1578.4166
$ws.Range("J137").Value = 5000.5
$ws.Range("K137").Value = 4735.2498
$ws.Range("L137").Value = 15001.5
$ws.Range("M137").Value = -2185.2498
$ws.Range("N137").Value = -20101.5

$ws.Range("H138").Value = 1844.2195
$ws.Range("I138").Value = 771.8929000000001
$ws.Range("J138").Value = 4153.846
$ws.Range("K138").Value = 2315.6787
$ws.Range("L138").Value = 12461.538
$ws.Range("M138").Value = 2824.3213
$ws.Range("N138").Value = -22741.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8445.788
$ws.Range("I32").Value = 2548.3125
$ws.Range("J32").Value = 26419.047
$ws.Range("K32").Value = 2548.3125
$ws.Range("L32").Value = 26419.047
$ws.Range("M32").Value = -2261.3125
$ws.Range("N32").Value = -26993.047

$ws.Range("H45").Value = 1963.5758
$ws.Range("I45").Value = 1917.1724
$ws.Range("J45").Value = 2300
$ws.Range("K45").Value = 1917.1724
$ws.Range("L45").Value = 2300
$ws.Range("M45").Value = -1540.1724
$ws.Range("N45").Value = -3054

$ws.Range("H74").Value = 8336511
$ws.Range("I74").Value = 11366468
$ws.Range("K74").Value = 11366468
$ws.Range("M74").Value = -11365594

$ws.Range("H77").Value = 8336511
$ws.Range("I77").Value = 11366468
$ws.Range("K77").Value = 56832340
$ws.Range("M77").Value = -56827972

$ws.Range("H132").Value = 3707.5715
$ws.Range("I132").Value = 3547.8333
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 10643.4999
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -8113.499899999999
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 863.7045000000001
$ws.Range("I134").Value = 530.7353000000001
$ws.Range("J134").Value = 1995.8
$ws.Range("K134").Value = 1592.2059
$ws.Range("L134").Value = 5987.4
$ws.Range("M134").Value = 942.7940999999998
$ws.Range("N134").Value = -11057.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2756.423
$ws.Range("I16").Value = 2628.05
$ws.Range("J16").Value = 3184.3333
$ws.Range("K16").Value = 2628.05
$ws.Range("L16").Value = 3184.3333
$ws.Range("M16").Value = -2341.05
$ws.Range("N16").Value = -3758.3333

$ws.Range("H31").Value = 20131.777
$ws.Range("I31").Value = 35533.137
$ws.Range("J31").Value = 2266.2
$ws.Range("K31").Value = 35533.137
$ws.Range("L31").Value = 2266.2
$ws.Range("M31").Value = -35238.137
$ws.Range("N31").Value = -2856.2

$ws.Range("H34").Value = 20131.777
$ws.Range("I34").Value = 35533.137
$ws.Range("J34").Value = 2266.2
$ws.Range("K34").Value = 35533.137
$ws.Range("L34").Value = 2266.2
$ws.Range("M34").Value = -35331.137
$ws.Range("N34").Value = -2670.2

$ws.Range("H58").Value = 837.7727
$ws.Range("I58").Value = 557.8684
$ws.Range("J58").Value = 1217.6428
$ws.Range("K58").Value = 557.8684
$ws.Range("L58").Value = 1217.6428
$ws.Range("M58").Value = -354.8684
$ws.Range("N58").Value = -1623.6428

$ws.Range("H113").Value = 2756.423
$ws.Range("I113").Value = 2628.05
$ws.Range("J113").Value = 3184.3333
$ws.Range("K113").Value = 2628.05
$ws.Range("L113").Value = 3184.3333
$ws.Range("M113").Value = -458.0500000000002
$ws.Range("N113").Value = -7524.3333

$ws.Range("H132").Value = 1975.3226
$ws.Range("I132").Value = 1311.375
$ws.Range("J132").Value = 4251.7144
$ws.Range("K132").Value = 3934.125
$ws.Range("L132").Value = 12755.1432
$ws.Range("M132").Value = -1404.125
$ws.Range("N132").Value = -17815.1432

$ws.Range("H136").Value = 837.7727
$ws.Range("I136").Value = 557.8684
$ws.Range("J136").Value = 1217.6428
$ws.Range("K136").Value = 1673.6052
$ws.Range("L136").Value = 3652.9284
$ws.Range("M136").Value = 876.3948
$ws.Range("N136").Value = -8752.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7937400
$ws.Range("I131").Value = 267.5
$ws.Range("J131").Value = 8772887
$ws.Range("K131").Value = 802.5
$ws.Range("L131").Value = 26318661
$ws.Range("M131").Value = 4237.5
$ws.Range("N131").Value = -26328741

$ws.Range("H132").Value = 882.1111
$ws.Range("I132").Value = 791.1818
$ws.Range("J132").Value = 1025
$ws.Range("K132").Value = 7120.6362
$ws.Range("L132").Value = 9225
$ws.Range("M132").Value = -4590.6362
$ws.Range("N132").Value = -14285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 688.1875
$ws.Range("I107").Value = 328.33334
$ws.Range("K107").Value = 328.33334
$ws.Range("M107").Value = 1591.66666

$ws.Range("H122").Value = 3745.5757
$ws.Range("I122").Value = 2958.647
$ws.Range("J122").Value = 4581.6875
$ws.Range("K122").Value = 8875.940999999999
$ws.Range("L122").Value = 13745.0625
$ws.Range("M122").Value = -6425.940999999999
$ws.Range("N122").Value = -18645.0625

$ws.Range("H126").Value = 2643.842
$ws.Range("I126").Value = 1648.6923
$ws.Range("J126").Value = 4800
$ws.Range("K126").Value = 4946.0769
$ws.Range("L126").Value = 14400
$ws.Range("M126").Value = -2476.0769
$ws.Range("N126").Value = -19340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -1388
$ws.Range("N7").Value = -1724

$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 958.6875
$ws.Range("I107").Value = 733.53845
$ws.Range("J107").Value = 1934.3334
$ws.Range("K107").Value = 2200.61535
$ws.Range("L107").Value = 5803.0002
$ws.Range("M107").Value = -280.61535
$ws.Range("N107").Value = -9643.0002

$ws.Range("H132").Value = 2217.8718
$ws.Range("I132").Value = 1551.4445
$ws.Range("J132").Value = 3717.3333
$ws.Range("K132").Value = 4654.333500000001
$ws.Range("L132").Value = 11151.9999
$ws.Range("M132").Value = -2124.333500000001
$ws.Range("N132").Value = -16211.9999

